$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.105.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.520.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.32"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.519.08"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.89"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.116.38"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.93"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.524.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.100.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.56%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.611"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.663.75"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.30"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.54%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.46"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.65%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "173.07"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -9.75%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.06"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.29"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.21%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.46"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.32%  "
